$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update values for rows 2-6 (2014/12 .. 2018/12) ---
# Row 2
$ws.Range("D2").Value = 652
$ws.Range("E2").Value = 77
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 86
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 66
$ws.Range("K2").Value = 874
$ws.Range("L2").Value = 225
$ws.Range("M2").Value = 649
$ws.Range("N2").Value = 649
$ws.Range("P2").Value = 162
$ws.Range("Q2").Value = 67
$ws.Range("R2").Value = -27
$ws.Range("S2").Value = -8
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 66
$ws.Range("V2").Value = 4
$ws.Range("W2").Value = 11.79
$ws.Range("X2").Value = 10.16
$ws.Range("Y2").Value = 10.65
$ws.Range("Z2").Value = 7.76
$ws.Range("AA2").Value = 34.62
$ws.Range("AB2").Value = 303.65
$ws.Range("AC2").Value = 204
$ws.Range("AD2").Value = 7.83
$ws.Range("AE2").Value = 2182
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = 45
$ws.Range("AH2").Value = 2.81
$ws.Range("AI2").Value = 20.2
$ws.Range("AJ2").Value = 32446151
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 627
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 75
$ws.Range("G3").Value = 79
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 63
$ws.Range("K3").Value = 920
$ws.Range("L3").Value = 213
$ws.Range("M3").Value = 706
$ws.Range("N3").Value = 706
$ws.Range("P3").Value = 162
$ws.Range("Q3").Value = 19
$ws.Range("R3").Value = 96
$ws.Range("S3").Value = -13
$ws.Range("T3").Value = 56
$ws.Range("U3").Value = -36
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 11.91
$ws.Range("X3").Value = 10.05
$ws.Range("Y3").Value = 9.300000000000001
$ws.Range("Z3").Value = 7.03
$ws.Range("AA3").Value = 30.21
$ws.Range("AB3").Value = 338
$ws.Range("AC3").Value = 194
$ws.Range("AD3").Value = 12.31
$ws.Range("AE3").Value = 2363
$ws.Range("AF3").Value = 1.01
$ws.Range("AG3").Value = 45
$ws.Range("AH3").Value = 1.88
$ws.Range("AI3").Value = 21.34
$ws.Range("AJ3").Value = 32446151
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 627
$ws.Range("E4").Value = 87
$ws.Range("F4").Value = 87
$ws.Range("G4").Value = 89
$ws.Range("H4").Value = 62
$ws.Range("I4").Value = 62
$ws.Range("K4").Value = 951
$ws.Range("L4").Value = 178
$ws.Range("M4").Value = 773
$ws.Range("N4").Value = 773
$ws.Range("P4").Value = 162
$ws.Range("Q4").Value = 57
$ws.Range("R4").Value = -59
$ws.Range("S4").Value = 36
$ws.Range("T4").Value = 27
$ws.Range("U4").Value = 30
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 13.84
$ws.Range("X4").Value = 9.84
$ws.Range("Y4").Value = 8.34
$ws.Range("Z4").Value = 6.59
$ws.Range("AA4").Value = 23.07
$ws.Range("AB4").Value = 373.56
$ws.Range("AC4").Value = 190
$ws.Range("AD4").Value = 15.65
$ws.Range("AE4").Value = 2432
$ws.Range("AF4").Value = 1.22
$ws.Range("AG4").Value = 65
$ws.Range("AH4").Value = 2.18
$ws.Range("AI4").Value = 33.49
$ws.Range("AJ4").Value = 32446151
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 562
$ws.Range("E5").Value = 62
$ws.Range("F5").Value = 62
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 47
$ws.Range("I5").Value = 47
$ws.Range("K5").Value = 903
$ws.Range("L5").Value = 122
$ws.Range("M5").Value = 781
$ws.Range("N5").Value = 781
$ws.Range("P5").Value = 162
$ws.Range("Q5").Value = -26
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = -18
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = -30
$ws.Range("V5").Value = 17
$ws.Range("W5").Value = 11.05
$ws.Range("X5").Value = 8.34
$ws.Range("Y5").Value = 6.03
$ws.Range("Z5").Value = 5.05
$ws.Range("AA5").Value = 15.6
$ws.Range("AB5").Value = 383.57
$ws.Range("AC5").Value = 144
$ws.Range("AD5").Value = 26.84
$ws.Range("AE5").Value = 2457
$ws.Range("AF5").Value = 1.58
$ws.Range("AG5").Value = 65
$ws.Range("AH5").Value = 1.68
$ws.Range("AI5").Value = 44.11
$ws.Range("AJ5").Value = 32446151
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 533
$ws.Range("E6").Value = 48
$ws.Range("F6").Value = 48
$ws.Range("G6").Value = 44
$ws.Range("H6").Value = 33
$ws.Range("I6").Value = 33
$ws.Range("K6").Value = 907
$ws.Range("L6").Value = 109
$ws.Range("M6").Value = 798
$ws.Range("N6").Value = 798
$ws.Range("P6").Value = 162
$ws.Range("Q6").Value = 28
$ws.Range("R6").Value = -34
$ws.Range("S6").Value = -20
$ws.Range("T6").Value = 19
$ws.Range("U6").Value = 9
$ws.Range("V6").Value = 18
$ws.Range("W6").Value = 9.02
$ws.Range("X6").Value = 6.28
$ws.Range("Y6").Value = 4.24
$ws.Range("Z6").Value = 3.7
$ws.Range("AA6").Value = 13.6
$ws.Range("AB6").Value = 394.27
$ws.Range("AC6").Value = 103
$ws.Range("AD6").Value = 20.71
$ws.Range("AE6").Value = 2511
$ws.Range("AF6").Value = 0.85
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 1.64
$ws.Range("AI6").Value = 33.25
$ws.Range("AJ6").Value = 32446151

# --- Clear forecast rows 7-9 data (keep A/B/C identifier cells) ---
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
